$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''57.597.64'
$ws.Range("E2").Value = '''  -1.68%  '

$ws.Range("D3").Value = '''2.434.82'
$ws.Range("E3").Value = '''  -2.17%  '

$ws.Range("E4").Value = '''  +0.02%  '

$ws.Range("D5").Value = '''513.70'
$ws.Range("E5").Value = '''  -2.63%  '

$ws.Range("D6").Value = '''130.09'
$ws.Range("E6").Value = '''  -3.06%  '

$ws.Range("E7").Value = '''  -0.11%  '

$ws.Range("D8").Value = '''0.549'
$ws.Range("E8").Value = '''  -2.09%  '

$ws.Range("D9").Value = '''2.444.76'
$ws.Range("E9").Value = '''  -1.99%  '

$ws.Range("E10").Value = '''  -0.35%  '

$ws.Range("D11").Value = '''0.0953'
$ws.Range("E11").Value = '''  -5.31%  '

$ws.Range("D12").Value = '''5.19'
$ws.Range("E12").Value = '''  -3.87%  '

$ws.Range("E13").Value = '''  -2.84%  '

$ws.Range("D14").Value = '''2.866.25'
$ws.Range("E14").Value = '''  -2.21%  '

$ws.Range("D15").Value = '''57.501.32'
$ws.Range("E15").Value = '''  -1.68%  '

$ws.Range("D16").Value = '''21.85'
$ws.Range("E16").Value = '''  -3.01%  '

$ws.Range("E17").Value = '''  -3.65%  '

$ws.Range("D18").Value = '''2.440.75'
$ws.Range("E18").Value = '''  -1.92%  '

$ws.Range("D19").Value = '''10.46'
$ws.Range("E19").Value = '''  -4.53%  '

$ws.Range("D20").Value = '''315.76'
$ws.Range("E20").Value = '''  -2.09%  '

$ws.Range("D21").Value = '''4.11'
$ws.Range("E21").Value = '''  -2.50%  '

$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '''  +0.10%  '

$ws.Range("D23").Value = '''5.68'
$ws.Range("E23").Value = '''  -2.62%  '

$ws.Range("D24").Value = '''63.46'
$ws.Range("E24").Value = '''  -1.62%  '

$ws.Range("D25").Value = '''0.408'
$ws.Range("E25").Value = '''  -1.70%  '

$ws.Range("D26").Value = '''0.996'
$ws.Range("E26").Value = '''  -0.25%  '

$ws.Range("E27").Value = '''  -1.93%  '

$ws.Range("D28").Value = '''7.20'
$ws.Range("E28").Value = '''  -3.60%  '

$ws.Range("D29").Value = '''170.13'
$ws.Range("E29").Value = '''  +1.76%  '

$ws.Range("D30").Value = '''0.0₃0723'
$ws.Range("E30").Value = '''  -4.30%  '

$ws.Range("D31").Value = '''6.26'
$ws.Range("E31").Value = '''  -3.31%  '

$ws.Range("D32").Value = '''1.67'
$ws.Range("E32").Value = '''  -3.40%  '

$ws.Range("E33").Value = '''  +2.46%  '

$ws.Range("D34").Value = '''0.998'
$ws.Range("E34").Value = '''  -0.04%  '

$ws.Range("D35").Value = '''0.997'
$ws.Range("E35").Value = '''  -0.17%  '

$ws.Range("D36").Value = '''17.71'
$ws.Range("E36").Value = '''  -3.22%  '

$ws.Range("D37").Value = '''1.29'
$ws.Range("E37").Value = '''  -4.15%  '

$ws.Range("D38").Value = '''3.91'
$ws.Range("E38").Value = '''  -2.27%  '

$ws.Range("D39").Value = '''36.20'
$ws.Range("E39").Value = '''  -1.40%  '

$ws.Range("D40").Value = '''1.46'
$ws.Range("E40").Value = '''  -2.92%  '

$ws.Range("D41").Value = '''0.781'
$ws.Range("E41").Value = '''  -2.38%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '''5.01'
$ws.Range("E42").Value = '''  +0.68%  '

$ws.Range("D43").Value = '''269.22'
$ws.Range("E43").Value = '''  -2.68%  '

$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '''3.38'
$ws.Range("E44").Value = '''  -4.57%  '

$ws.Range("D45").Value = '''0.589'
$ws.Range("E45").Value = '''  -1.72%  '

$ws.Range("D46").Value = '''0.0907'
$ws.Range("E46").Value = '''  -0.94%  '

$ws.Range("D47").Value = '''120.37'
$ws.Range("E47").Value = '''  -5.94%  '

$ws.Range("D48").Value = '''0.0485'
$ws.Range("E48").Value = '''  -2.56%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '''0.0211'
$ws.Range("E49").Value = '''  -2.77%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''16.59'
$ws.Range("E50").Value = '''  -3.76%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '''1.706.50'
$ws.Range("E51").Value = '''  -1.86%  '
